$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("A15").Value = "Osmar Ali Baba"
$ws.Range("A10").Value = "Pedro Ramoneda Franquista"
$ws.Range("A11").Value = "Joaquin Fernández OfensiveMan"
$ws.Range("A14").Value = "Alejandro Francés Rubio"
$ws.Range("A16").Value = "Diego Santolaya Firulais"

$ws.Columns.Item(1).ColumnWidth = 29.5703125

$ws.Range("A16").Select()
